$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1245411268104184
$ws.Range("E2").Value = 8.744009178710868
$ws.Range("F2").Value = 24.24618101016198
